# Update example model specification sheets:
# add a "type" column (valued "static") to the costs and utilities sheets,
# inserted right after the "state" column, and make "costs" the active sheet.

$wb = $excel.ActiveWorkbook

# --- costs sheet -----------------------------------------------------------
$wsCosts = $wb.Worksheets.Item("costs")
$wsCosts.Range("B1").EntireColumn.Insert() | Out-Null
$wsCosts.Range("B1").Value = "type"
$wsCosts.Range("B2:B4").Value = "static"

# --- utilities sheet ---------------------------------------------------------
$wsUtil = $wb.Worksheets.Item("utilities")
$wsUtil.Range("B1").EntireColumn.Insert() | Out-Null
$wsUtil.Range("B1").Value = "type"
$wsUtil.Range("B2:B4").Value = "static"
$wsUtil.Range("B2:B4").Select() | Out-Null

# --- make "costs" the active/selected sheet, matching the new selection ----
$wsCosts.Select() | Out-Null
$wsCosts.Range("B2:B4").Select() | Out-Null
